# update error message handling
# Row 8 ("Error message" / "Sprint 4") moves from status "Open" to "Done":
#  - D8 value changes from "Open" to "Done"
#  - the whole row (A8:E8) is re-styled to match the other "Done" rows:
#      green fill (same RGB as used elsewhere: 5AC664) on A8:E8
#      left-aligned text on A8, B8, C8, D8
#      wrapped text on C8 (long description column)
#      E8 (Assignee) keeps default (general) alignment/no-wrap, just gets the fill
#  - the active selection moves from D2 to C5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status value: Open -> Done ---
$ws.Range("D8").Value = "Done"

# --- Re-style row 8 like the other "Done" rows (green fill + alignment) ---
$green = 6604378  # BGR encoding of RGB 5AC664 (the fill already used by the other "Done" rows)

$ws.Range("A8:E8").Interior.Color = $green

$ws.Range("A8:B8").HorizontalAlignment = -4131   # xlLeft
$ws.Range("D8").HorizontalAlignment = -4131      # xlLeft

$ws.Range("C8").HorizontalAlignment = -4131      # xlLeft
$ws.Range("C8").WrapText = $true

# --- Update the active selection/active cell ---
[void]$ws.Range("C5").Select()
